$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = 0.1454560131662833
$ws.Range("D3").Value  = 0.2960932320440638
$ws.Range("D4").Value  = 0.4139825968325628
$ws.Range("D5").Value  = 0.6030192320592257
$ws.Range("D6").Value  = 0.8080822157080472
$ws.Range("D7").Value  = 1.184890562614349
$ws.Range("D8").Value  = 0.222903107808925
$ws.Range("D9").Value  = 0.5030254831398444
$ws.Range("D10").Value = 0.6948495953525461
$ws.Range("D11").Value = 1.085632739753365
$ws.Range("D12").Value = 1.452731959033275
$ws.Range("D13").Value = 2.49027611399768
